{"js": "const newValues = [\n  [\"37+25=\", \"53+1=\", \"9+49=\", \"0+76=\", \"6+1=\"],\n  [\"94-84=\", \"20-2=\", \"99-96=\", \"45+42=\", \"23+3=\"],\n  [\"7-6=\", \"3+18=\", \"71-5=\", \"67+20=\", \"46+18=\"],\n  [\"5+7=\", \"90-55=\", \"59+39=\", \"79-2=\", \"16+11=\"],\n  [\"68-61=\", \"86-71=\", \"34+18=\", \"36+33=\", \"56-53=\"],\n  [\"37-32=\", \"73-24=\", \"88-31=\", \"6+45=\", \"64-28=\"],\n  [\"4+74=\", \"63+20=\", \"80-34=\", \"89-4=\", \"47+37=\"],\n  [\"36+31=\", \"52+40=\", \"99-89=\", \"41+6=\", \"64-17=\"],\n  [\"4-2=\", \"34+20=\", \"26+34=\", \"97-24=\", \"63-34=\"],\n  [\"25+4=\", \"49-35=\", \"6+53=\", \"51+21=\", \"1+49=\"],\n  [\"13+83=\", \"50-15=\", \"23-18=\", \"55+28=\", \"8+23=\"],\n  [\"90-40=\", \"92-49=\", \"72-59=\", \"34+41=\", \"4+37=\"],\n  [\"51+6=\", \"35-20=\", \"30+61=\", \"75+17=\", \"57-54=\"],\n  [\"8+90=\", \"89-57=\", \"16+76=\", \"4+61=\", \"88-64=\"],\n  [\"83-33=\", \"35-34=\", \"7+77=\", \"46-1=\", \"78-50=\"],\n  [\"27+40=\", \"40-24=\", \"92-41=\", \"80-2=\", \"52-13=\"],\n  [\"49+2=\", \"16+43=\", \"60+22=\", \"7+64=\", \"10+71=\"],\n  [\"30-14=\", \"67+25=\", \"7+77=\", \"81-59=\", \"53-52=\"],\n  [\"87-47=\", \"34-25=\", \"37+11=\", \"95-5=\", \"44+55=\"],\n  [\"48+6=\", \"12+81=\", \"0+71=\", \"98-55=\", \"71-32=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"37+25=\",\"53+1=\",\"9+49=\",\"0+76=\",\"6+1=\"),\n    @(\"94-84=\",\"20-2=\",\"99-96=\",\"45+42=\",\"23+3=\"),\n    @(\"7-6=\",\"3+18=\",\"71-5=\",\"67+20=\",\"46+18=\"),\n    @(\"5+7=\",\"90-55=\",\"59+39=\",\"79-2=\",\"16+11=\"),\n    @(\"68-61=\",\"86-71=\",\"34+18=\",\"36+33=\",\"56-53=\"),\n    @(\"37-32=\",\"73-24=\",\"88-31=\",\"6+45=\",\"64-28=\"),\n    @(\"4+74=\",\"63+20=\",\"80-34=\",\"89-4=\",\"47+37=\"),\n    @(\"36+31=\",\"52+40=\",\"99-89=\",\"41+6=\",\"64-17=\"),\n    @(\"4-2=\",\"34+20=\",\"26+34=\",\"97-24=\",\"63-34=\"),\n    @(\"25+4=\",\"49-35=\",\"6+53=\",\"51+21=\",\"1+49=\"),\n    @(\"13+83=\",\"50-15=\",\"23-18=\",\"55+28=\",\"8+23=\"),\n    @(\"90-40=\",\"92-49=\",\"72-59=\",\"34+41=\",\"4+37=\"),\n    @(\"51+6=\",\"35-20=\",\"30+61=\",\"75+17=\",\"57-54=\"),\n    @(\"8+90=\",\"89-57=\",\"16+76=\",\"4+61=\",\"88-64=\"),\n    @(\"83-33=\",\"35-34=\",\"7+77=\",\"46-1=\",\"78-50=\"),\n    @(\"27+40=\",\"40-24=\",\"92-41=\",\"80-2=\",\"52-13=\"),\n    @(\"49+2=\",\"16+43=\",\"60+22=\",\"7+64=\",\"10+71=\"),\n    @(\"30-14=\",\"67+25=\",\"7+77=\",\"81-59=\",\"53-52=\"),\n    @(\"87-47=\",\"34-25=\",\"37+11=\",\"95-5=\",\"44+55=\"),\n    @(\"48+6=\",\"12+81=\",\"0+71=\",\"98-55=\",\"71-32=\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
